$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 13383.565
$ws.Range("J17").Value = 13383.565
$ws.Range("L17").Value = 40150.695
$ws.Range("N17").Value = -40486.695
$ws.Range("H33").Value = 2847.3845
$ws.Range("I33").Value = 1882.1
$ws.Range("K33").Value = 1882.1
$ws.Range("M33").Value = -1653.1
$ws.Range("H74").Value = 4874.364
$ws.Range("I74").Value = 4162.25
$ws.Range("K74").Value = 4162.25
$ws.Range("M74").Value = -3226.25
$ws.Range("H77").Value = 4874.364
$ws.Range("I77").Value = 4162.25
$ws.Range("K77").Value = 20811.25
$ws.Range("M77").Value = -16131.25
$ws.Range("H92").Value = 800.2857
$ws.Range("I92").Value = 767
$ws.Range("K92").Value = 767
$ws.Range("M92").Value = 481
$ws.Range("H97").Value = 2000
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H100").Value = 2725.9375
$ws.Range("I100").Value = 2699.7778
$ws.Range("K100").Value = 2699.7778
$ws.Range("M100").Value = -2158.7778
$ws.Range("H103").Value = 31250388
$ws.Range("I103").Value = 360
$ws.Range("J103").Value = 45454944
$ws.Range("K103").Value = 1080
$ws.Range("L103").Value = 136364832
$ws.Range("M103").Value = -494
$ws.Range("N103").Value = -136366004
$ws.Range("H138").Value = 1742.305
$ws.Range("J138").Value = 2134.5833
$ws.Range("L138").Value = 6403.749899999999
$ws.Range("N138").Value = -16683.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 957.125
$ws.Range("I61").Value = 957.125
$ws.Range("K61").Value = 957.125
$ws.Range("M61").Value = -745.125
$ws.Range("H63").Value = 2364.7693
$ws.Range("I63").Value = 2367.5454
$ws.Range("J63").Value = 2349.5
$ws.Range("K63").Value = 2367.5454
$ws.Range("L63").Value = 2349.5
$ws.Range("M63").Value = -1681.5454
$ws.Range("N63").Value = -3721.5
$ws.Range("H66").Value = 2364.7693
$ws.Range("I66").Value = 2367.5454
$ws.Range("J66").Value = 2349.5
$ws.Range("K66").Value = 11837.727
$ws.Range("L66").Value = 11747.5
$ws.Range("M66").Value = -8405.726999999999
$ws.Range("N66").Value = -18611.5
$ws.Range("H74").Value = 3022.625
$ws.Range("I74").Value = 2670.182
$ws.Range("J74").Value = 3798
$ws.Range("K74").Value = 2670.182
$ws.Range("L74").Value = 3798
$ws.Range("M74").Value = -1796.182
$ws.Range("N74").Value = -5546
$ws.Range("H77").Value = 3022.625
$ws.Range("I77").Value = 2670.182
$ws.Range("J77").Value = 3798
$ws.Range("K77").Value = 13350.91
$ws.Range("L77").Value = 18990
$ws.Range("M77").Value = -8982.91
$ws.Range("N77").Value = -27726
$ws.Range("H122").Value = 3056
$ws.Range("I122").Value = 3056
$ws.Range("K122").Value = 9168
$ws.Range("M122").Value = -6718
$ws.Range("H132").Value = 11226.333
$ws.Range("I132").Value = 13608.9375
$ws.Range("K132").Value = 40826.8125
$ws.Range("M132").Value = -38296.8125
$ws.Range("H136").Value = 957.125
$ws.Range("I136").Value = 957.125
$ws.Range("K136").Value = 2871.375
$ws.Range("M136").Value = -321.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H130").Value = 75624.875
$ws.Range("J130").Value = 75624.875
$ws.Range("L130").Value = 75624.875
$ws.Range("N130").Value = -85664.875
$ws.Range("H134").Value = 1924.174
$ws.Range("I134").Value = 1463.7222
$ws.Range("K134").Value = 4391.1666
$ws.Range("M134").Value = -1856.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 35714570
$ws.Range("J7").Value = 407.46155
$ws.Range("L7").Value = 407.46155
$ws.Range("N7").Value = -633.46155
$ws.Range("H59").Value = 18332.666
$ws.Range("J59").Value = 18332.666
$ws.Range("L59").Value = 18332.666
$ws.Range("N59").Value = -20622.666
$ws.Range("H132").Value = 2013.3846
$ws.Range("I132").Value = 2088.4783
$ws.Range("J132").Value = 1437.6666
$ws.Range("K132").Value = 6265.4349
$ws.Range("L132").Value = 4312.9998
$ws.Range("M132").Value = -3735.4349
$ws.Range("N132").Value = -9372.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42829020
$ws.Range("I4").Value = 3381666.8
$ws.Range("K4").Value = 10145000.4
$ws.Range("M4").Value = -10144888.4
$ws.Range("H38").Value = 105.3
$ws.Range("I38").Value = 59
$ws.Range("K38").Value = 177
$ws.Range("M38").Value = 170
$ws.Range("H98").Value = 1666.1666
$ws.Range("I98").Value = 1334.3334
$ws.Range("J98").Value = 1998
$ws.Range("K98").Value = 4003.0002
$ws.Range("L98").Value = 5994
$ws.Range("M98").Value = -2505.0002
$ws.Range("N98").Value = -8990
$ws.Range("H136").Value = 3948.625
$ws.Range("I136").Value = 1994
$ws.Range("K136").Value = 5982
$ws.Range("M136").Value = -882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9999.5
$ws.Range("J92").Value = 9999.5
$ws.Range("L92").Value = 9999.5
$ws.Range("N92").Value = -13743.5
$ws.Range("H104").Value = 20223.334
$ws.Range("J104").Value = 20223.334
$ws.Range("L104").Value = 20223.334
$ws.Range("N104").Value = -27211.334
$ws.Range("H126").Value = 3957
$ws.Range("J126").Value = 3914
$ws.Range("L126").Value = 11742
$ws.Range("N126").Value = -16682
$ws.Range("H134").Value = 38884
$ws.Range("J134").Value = 38884
$ws.Range("L134").Value = 116652
$ws.Range("N134").Value = -121722

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8862.375
$ws.Range("I40").Value = 9484
$ws.Range("K40").Value = 9484
$ws.Range("M40").Value = -9348
$ws.Range("H55").Value = 411.1111
$ws.Range("I55").Value = 437.85715
$ws.Range("J55").Value = 317.5
$ws.Range("K55").Value = 437.85715
$ws.Range("L55").Value = 317.5
$ws.Range("M55").Value = -264.85715
$ws.Range("N55").Value = -663.5
$ws.Range("H136").Value = 2915.6785
$ws.Range("I136").Value = 2635.2856
$ws.Range("J136").Value = 3196.0715
$ws.Range("K136").Value = 7905.8568
$ws.Range("L136").Value = 9588.2145
$ws.Range("M136").Value = -5355.8568
$ws.Range("N136").Value = -14688.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2226.0356
$ws.Range("I122").Value = 1798.8182
$ws.Range("K122").Value = 5396.4546
$ws.Range("M122").Value = -2946.4546
$ws.Range("H126").Value = 3205.3333
$ws.Range("I126").Value = 2783.8572
$ws.Range("K126").Value = 8351.571599999999
$ws.Range("M126").Value = -5881.571599999999
$ws.Range("H132").Value = 15724.5
$ws.Range("I132").Value = 28304
$ws.Range("J132").Value = 7862.3125
$ws.Range("K132").Value = 84912
$ws.Range("L132").Value = 23586.9375
$ws.Range("M132").Value = -82382
$ws.Range("N132").Value = -28646.9375
$ws.Range("H136").Value = 893.5833
$ws.Range("I136").Value = 892.6667
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 2678.0001
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -128.0001000000002
$ws.Range("N136").Value = -7800
